# Applies the "Updated cryptos list" data refresh to Sheet1.
# Most text cells in columns D (Price) and E (Volume 1h) look like plain
# numbers/percentages but are stored as literal text (inline strings) in
# the workbook. Writing a numeric-looking string straight into .Value
# lets Excel auto-convert it to a real number, which would corrupt the
# original text formatting (e.g. "45.85" -> 45.85 as Double). To avoid
# that we force the cell to Text format ("@") before assigning the
# value, then restore the cell style to Normal so we don't leave stray
# formatting behind on cells that were plain (unstyled) text before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "43.774.39"
Set-TextValue "E2" "  -0.07%  "
Set-TextValue "D3" "2.291.89"
Set-TextValue "E3" "  -1.10%  "
Set-TextValue "E4" "  -0.15%  "
Set-TextValue "D5" "103.43"
Set-TextValue "E5" "  +6.48%  "
Set-TextValue "D6" "271.13"
Set-TextValue "E6" "  -0.48%  "
Set-TextValue "E7" "  -0.52%  "
Set-TextValue "E8" "  +0.09%  "
Set-TextValue "D9" "0.606"
Set-TextValue "E9" "  -3.29%  "
Set-TextValue "D10" "45.85"
Set-TextValue "E10" "  +1.05%  "
Set-TextValue "E11" "  -1.70%  "
Set-TextValue "D12" "8.23"
Set-TextValue "E12" "  +2.47%  "
Set-TextValue "E13" "  +1.71%  "
Set-TextValue "E14" "  -0.05%  "
Set-TextValue "E15" "  -2.31%  "
Set-TextValue "D16" "2.293.54"
Set-TextValue "E16" "  -1.21%  "
Set-TextValue "D17" "43.788.18"
Set-TextValue "E17" "  +0.04%  "
Set-TextValue "E18" "  +0.18%  "
Set-TextValue "D19" "6.28"
Set-TextValue "E19" "  -1.91%  "
Set-TextValue "D20" "72.25"
Set-TextValue "E21" "  +10.35%  "
Set-TextValue "D22" "233.57"
Set-TextValue "E22" "  -2.59%  "
Set-TextValue "D23" "2.93"
Set-TextValue "E23" "  +14.99%  "
Set-TextValue "D24" "9.29"
Set-TextValue "E24" "  -1.39%  "
Set-TextValue "E25" "  +0.04%  "
Set-TextValue "D26" "11.30"
Set-TextValue "E26" "  -0.60%  "
Set-TextValue "B27" "WEMIXToken"
Set-TextValue "C27" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D27" "3.45"
Set-TextValue "E27" "  -1.62%  "
Set-TextValue "B28" "InjectiveProtocol"
Set-TextValue "C28" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D28" "40.29"
Set-TextValue "E28" "  +5.21%  "
Set-TextValue "E29" "  -3.10%  "
Set-TextValue "D30" "177.69"
Set-TextValue "E30" "  +1.47%  "
Set-TextValue "D31" "21.83"
Set-TextValue "E31" "  -2.60%  "
Set-TextValue "D32" "0.0901"
Set-TextValue "E32" "  -1.90%  "
Set-TextValue "D33" "5.53"
Set-TextValue "E33" "  +0.55%  "
Set-TextValue "E34" "  +10.04%  "
Set-TextValue "E35" "  -0.35%  "
Set-TextValue "E36" "  +1.68%  "
Set-TextValue "E37" "  -1.48%  "
Set-TextValue "E38" "  +4.94%  "
Set-TextValue "B39" "Algorand"
Set-TextValue "C39" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D39" "0.237"
Set-TextValue "E39" "  -3.39%  "
Set-TextValue "B40" "LidoDAOToken"
Set-TextValue "C40" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D40" "2.33"
Set-TextValue "E40" "  -2.92%  "
Set-TextValue "E41" "  -2.86%  "
Set-TextValue "D42" "65.44"
Set-TextValue "E42" "  +4.34%  "
Set-TextValue "D43" "12.16"
Set-TextValue "E43" "  -1.32%  "
Set-TextValue "D44" "5.34"
Set-TextValue "E44" "  -0.37%  "
Set-TextValue "D45" "8.79"
Set-TextValue "E45" "  -4.70%  "
Set-TextValue "E46" "  -1.52%  "
Set-TextValue "E47" "  +2.21%  "
Set-TextValue "D48" "99.33"
Set-TextValue "E48" "  -1.14%  "
Set-TextValue "E49" "  +12.21%  "
Set-TextValue "D50" "0.439"
Set-TextValue "E50" "  +4.29%  "
Set-TextValue "D51" "2.526.97"
Set-TextValue "E51" "  -0.67%  "
